$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 495
    $ws.Range("F3").Value = 3364
    $ws.Range("F4").Value = 91
}
